$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.108.38"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "3.562.17"
$ws.Range("E3").Value = "  +4.56%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'599.42"
$ws.Range("E5").Value = "  +3.40%  "

$ws.Range("D6").Value = "'137.86"
$ws.Range("E6").Value = "  +3.64%  "

$ws.Range("D7").Value = "3.560.54"
$ws.Range("E7").Value = "  +4.53%  "

$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +3.59%  "

$ws.Range("E10").Value = "  +3.53%  "

$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").Value = "'0.388"
$ws.Range("E12").Value = "  +4.69%  "

$ws.Range("D13").Value = "4.163.21"
$ws.Range("E13").Value = "  +4.36%  "

$ws.Range("E14").Value = "  +4.07%  "

$ws.Range("D15").Value = "'27.35"
$ws.Range("E15").Value = "  +5.72%  "

$ws.Range("D16").Value = "3.565.17"
$ws.Range("E16").Value = "  +4.05%  "

$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D18").Value = "64.926.87"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").Value = "'10.15"
$ws.Range("E19").Value = "  +8.21%  "

$ws.Range("D20").Value = "'14.42"
$ws.Range("E20").Value = "  +7.64%  "

$ws.Range("E21").Value = "  +2.89%  "

$ws.Range("D22").Value = "'391.69"
$ws.Range("E22").Value = "  +3.44%  "

$ws.Range("E23").Value = "  +8.15%  "

$ws.Range("D24").Value = "3.703.91"
$ws.Range("E24").Value = "  +4.35%  "

$ws.Range("D25").Value = "'74.20"
$ws.Range("E25").Value = "  +3.98%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "'0.0000118"
$ws.Range("E27").Value = "  +14.58%  "

$ws.Range("D28").Value = "'7.71"
$ws.Range("E28").Value = "  +8.66%  "

$ws.Range("E29").Value = "  +0.38%  "

$ws.Range("D30").Value = "'2.30"
$ws.Range("E30").Value = "  +5.93%  "

$ws.Range("D31").Value = "'8.33"
$ws.Range("E31").Value = "  +5.83%  "

$ws.Range("D32").Value = "3.567.09"
$ws.Range("E32").Value = "  +4.05%  "

$ws.Range("E33").Value = "  +24.90%  "

$ws.Range("D34").Value = "'24.06"
$ws.Range("E34").Value = "  +5.95%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +2.33%  "

$ws.Range("D37").Value = "'170.17"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").Value = "'6.92"
$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  +8.71%  "

$ws.Range("D40").Value = "'5.02"
$ws.Range("E40").Value = "  +11.88%  "

$ws.Range("E41").Value = "  +7.86%  "

$ws.Range("D42").Value = "'0.828"
$ws.Range("E42").Value = "  +3.52%  "

$ws.Range("D43").Value = "'26.82"
$ws.Range("E43").Value = "  +21.04%  "

$ws.Range("D44").Value = "'42.65"
$ws.Range("E44").Value = "  +1.97%  "

$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("E46").Value = "  +10.82%  "

$ws.Range("D47").Value = "'4.47"
$ws.Range("E47").Value = "  +6.36%  "

$ws.Range("E48").Value = "  +5.23%  "

$ws.Range("D49").Value = "2.499.85"
$ws.Range("E49").Value = "  +14.13%  "

$ws.Range("E50").Value = "  +7.76%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'303.39"
$ws.Range("E51").Value = "  +11.48%  "
